$wb = $excel.ActiveWorkbook

# Work on the 'master-flag-group' sheet: fill in the "Group" column (C)
# for the rows that were still blank, finishing the grouping work for
# the Mesos Master configuration flags.
$ws = $wb.Worksheets.Item("master-flag-group")

$ws.Range("C25").Value = "core"
$ws.Range("C26").Value = "security"
$ws.Range("C27").Value = "core"
$ws.Range("C28").Value = "allocation"
$ws.Range("C29").Value = "core"
$ws.Range("C30").Value = "core"
$ws.Range("C31").Value = "allocation"

# Make this the active sheet/selection, matching the saved view state.
$ws.Activate()
$ws.Range("C32").Select()
